# Weekly data refresh: two new weekly price records are inserted into the
# "Espinaca" (spinach) sheet, pushing the existing rows below them down.
#
# New record #1 is inserted at row 25 (everything that was at row 25..147
# shifts down to 26..148).
# New record #2 is inserted at row 72 of the *resulting* sheet (everything
# that was at row 72..148 at that point shifts down to 73..149).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new weekly record #1 at row 25 -------------------------------
$ws.Rows.Item(25).Insert()

$row1 = New-Object 'object[,]' 1,18
$row1[0,0]  = 11
$row1[0,1]  = "Vega Monumental Concepción"
$row1[0,2]  = "Bíobío"
$row1[0,3]  = [datetime]"2022-11-09"
$row1[0,4]  = 8
$row1[0,5]  = 100112012
$row1[0,6]  = "Espinaca"
$row1[0,7]  = "Sin especificar"
$row1[0,8]  = "Primera"
$row1[0,9]  = 130
$row1[0,10] = 14000
$row1[0,11] = 15000
$row1[0,12] = 14615
$row1[0,13] = "`$/cuna 10 kilos"
$row1[0,14] = "Provincia de Chacabuco"
$row1[0,15] = 1462
$row1[0,16] = 10
$row1[0,17] = "Hortaliza"
$ws.Range("A25:R25").Value = $row1

# --- Insert new weekly record #2 at row 72 (post first insert) -----------
$ws.Rows.Item(72).Insert()

$row2 = New-Object 'object[,]' 1,18
$row2[0,0]  = 11
$row2[0,1]  = "Vega Monumental Concepción"
$row2[0,2]  = "Bíobío"
$row2[0,3]  = [datetime]"2022-02-24"
$row2[0,4]  = 8
$row2[0,5]  = 100112012
$row2[0,6]  = "Espinaca"
$row2[0,7]  = "Sin especificar"
$row2[0,8]  = "Primera"
$row2[0,9]  = 60
$row2[0,10] = 14000
$row2[0,11] = 15000
$row2[0,12] = 14500
$row2[0,13] = "`$/cuna 10 kilos"
$row2[0,14] = "Región Metropolitana"
$row2[0,15] = 1450
$row2[0,16] = 10
$row2[0,17] = "Hortaliza"
$ws.Range("A72:R72").Value = $row2
